# Corrected from Simon comments
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Swap the values in C10 and C11 (dx @ center <-> dx @ 1 mm FOV)
$ws.Range("C10").Value = 10
$ws.Range("C11").Value = 128

# Move the active selection on the sheet to C12
$ws.Range("C12").Select()

# Reflect the updated window/view position (best effort; engine may keep
# the original xWindow/yWindow values as they are not separately modeled)
$excel.ActiveWindow.Left = 1860
